# "changed test case 5"
# Test case 5 (row 6, the 5th data row after the INPUT/OUTPUT header row)
# had its input expression changed from "a ^ (b ^ (c ^ d))" to "((a^b)^c)^d".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "((a^b)^c)^d"

# Reflect the user's selection landing on the edited cell after the change.
$ws.Range("A6").Select()
